# Project DesignFirst / Main.xlsx - save
# Update the "Integer max" value for rule R20 (row 10) on the Rules sheet
# from 18 to 100.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 100

$wb.Save()
